$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Date" column (B) entirely -- values shift left:
# old C (Budgeted) -> new B, old D (Spent) -> new C, old E (Remaining formula) -> new D
$ws.Columns("B:B").Delete()

# Correct the "Spent" values that had been entered in the wrong category:
# Rent spent 26, Gas spent 43 (category values now live in column C)
$ws.Range("C2").Value = 26
$ws.Range("C3").Value = 43
